$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): existing "Alt" labels split into "mean"/"std" pairs and
# three new Alt windows (Three/Five/Ten Year) get their own mean/std columns.
# ---------------------------------------------------------------------------
$headers = @(
    "Algorithm",
    "One Year Alt mean", "One Year Alt std",
    "Two Year Alt mean", "Two Year Alt std",
    "Three Year Alt mean", "Three Year Alt std",
    "Five Year Alt mean", "Five Year Alt std",
    "Ten Year Alt mean", "Ten Year Alt std"
)
$cols = @("B","C","D","E","F","G","H","I","J","K","L")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# New header cells (H1:L1) need the same bold/centered/bordered style as the
# existing header cells -- copy formatting only (not value) from B1.
$ws.Range("B1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)


# ---------------------------------------------------------------------------
# Data rows (2-8): algorithm name + 5 mean/std pairs per row.
# Row that used to hold "NB" (old row 9) is removed entirely.
# "CART" is renamed to "DTREE".
# ---------------------------------------------------------------------------
$data = @(
    @("LR",    0.8537243401759531, 0.02182000774888037, 0.8441947049089906, 0.03101234304715023, 0.8297123164519029, 0.02327315754081393, 0.813390261627907,  0.04873765251553298, 0.8277777777777778, 0.03055240666144628),
    @("LDA",   0.8355969836614998, 0.02261034016097998, 0.8191165655451371, 0.03166668956409304, 0.8070971930876036, 0.0229572952456023,  0.7962996608527131, 0.05652256670606583, 0.8095444238301381, 0.04296341918059716),
    @("KNN",   0.8867239212400504, 0.02168690359069748, 0.8929306857878286, 0.01799549314157119, 0.8996703626011389, 0.02177180518798494, 0.8927507267441861, 0.03062021073755635, 0.8874871160585446, 0.03597734512467365),
    @("DTREE", 0.7773732718894009, 0.03558298844132256, 0.7682800147085861, 0.02761666988297672, 0.7611227649585456, 0.02681576108957025, 0.7566133720930233, 0.03418277887336316, 0.7284786641929499, 0.04971540871817492),
    @("RTREE", 0.880242982823628,  0.02250277223500797, 0.8665241772384631, 0.0335667758566429,  0.8784786734591948, 0.02383846342041215, 0.8756116763565892, 0.04326656056172588, 0.8794681508967223, 0.02307238220686191),
    @("XTREE", 0.8608378718056138, 0.02387442014759547, 0.8577128148556719, 0.03117809915692107, 0.8587054240335631, 0.02903352372579127, 0.8414304748062016, 0.04559338484830374, 0.8449391877963306, 0.02214432309171755),
    @("SVM",   0.8737704231252618, 0.02020405379965951, 0.861803640375069,  0.03004434341382434, 0.8678953151533314, 0.02906690953958588, 0.8779796511627908, 0.03535374557925629, 0.8662440733869305, 0.04310377143565268)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $ws.Range("A" + $rowNum).Value = $r
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value = $data[$r][$c]
    }
}

# ---------------------------------------------------------------------------
# Drop the old trailing row (previously row 9, the "NB" entry) which no
# longer has data in the new layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Delete()
